# Append two new daily rows (2025-11-30) for the two stations, mirroring the
# formatting/layout of the existing rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (59) into the two new rows
# (60 and 61) before writing values into them.
$ws.Range("A59:F59").Copy() | Out-Null
$ws.Range("A60:F60").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A59:F59").Copy() | Out-Null
$ws.Range("A61:F61").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Row 60 - 四方坪站 (station 1) for date serial 45991 (2025-11-30)
$ws.Cells.Item(60, 1).Value = 45991
$ws.Cells.Item(60, 2).Value = "四方坪站"
$ws.Cells.Item(60, 3).Value = 8853.2999999999993
$ws.Cells.Item(60, 4).Value = 8028.17
$ws.Cells.Item(60, 5).Value = 2831.83
$ws.Cells.Item(60, 6).Value = 372

# Row 61 - 高岭站 (station 2) for date serial 45991 (2025-11-30)
$ws.Cells.Item(61, 1).Value = 45991
$ws.Cells.Item(61, 2).Value = "高岭站"
$ws.Cells.Item(61, 3).Value = 5103.96
$ws.Cells.Item(61, 4).Value = 4370.12
$ws.Cells.Item(61, 5).Value = 1372.71
$ws.Cells.Item(61, 6).Value = 176

# Update selection to mirror the saved view state (I59 instead of G59)
$ws.Range("I59").Select() | Out-Null
